# Updates betting-odds cells F2:AO4 on the active sheet to the new values
# captured in the 2025-12-13 Betfair Back/Lay export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: assigning arrays directly to a multi-cell Range.Value does not
# reliably persist in this runtime, so each cell is set individually
# via Cells.Item(row, column).

# Row 2
$ws.Cells.Item(2, 6).Value = 1.08  # F2
$ws.Cells.Item(2, 7).Value = 1.14  # G2
$ws.Cells.Item(2, 8).Value = 220  # H2
$ws.Cells.Item(2, 9).Value = 1000  # I2
$ws.Cells.Item(2, 10).Value = 1.01  # J2
$ws.Cells.Item(2, 11).Value = 14  # K2
$ws.Cells.Item(2, 12).Value = 0  # L2
$ws.Cells.Item(2, 13).Value = 0  # M2
$ws.Cells.Item(2, 14).Value = 1.02  # N2
$ws.Cells.Item(2, 15).Value = 1.08  # O2
$ws.Cells.Item(2, 16).Value = 1.01  # P2
$ws.Cells.Item(2, 17).Value = 1.68  # Q2
$ws.Cells.Item(2, 18).Value = 1.01  # R2
$ws.Cells.Item(2, 19).Value = 100  # S2
$ws.Cells.Item(2, 20).Value = 1.45  # T2
$ws.Cells.Item(2, 21).Value = 1.01  # U2
$ws.Cells.Item(2, 22).Value = 1.02  # V2
$ws.Cells.Item(2, 23).Value = 1.03  # W2
$ws.Cells.Item(2, 24).Value = 1000  # X2
$ws.Cells.Item(2, 25).Value = 1000  # Y2
$ws.Cells.Item(2, 26).Value = 1000  # Z2
$ws.Cells.Item(2, 27).Value = 1000  # AA2
$ws.Cells.Item(2, 28).Value = 1.3  # AB2
$ws.Cells.Item(2, 29).Value = 16.5  # AC2
$ws.Cells.Item(2, 30).Value = 1000  # AD2
$ws.Cells.Item(2, 31).Value = 1000  # AE2
$ws.Cells.Item(2, 32).Value = 11.5  # AF2
$ws.Cells.Item(2, 33).Value = 1000  # AG2
$ws.Cells.Item(2, 34).Value = 1000  # AH2
$ws.Cells.Item(2, 35).Value = 1000  # AI2
$ws.Cells.Item(2, 36).Value = 1000  # AJ2
$ws.Cells.Item(2, 37).Value = 1000  # AK2
$ws.Cells.Item(2, 38).Value = 1000  # AL2
$ws.Cells.Item(2, 39).Value = 1000  # AM2
$ws.Cells.Item(2, 40).Value = 1000  # AN2
$ws.Cells.Item(2, 41).Value = 1000  # AO2

# Row 3
$ws.Cells.Item(3, 6).Value = 2.72  # F3
$ws.Cells.Item(3, 7).Value = 2.94  # G3
$ws.Cells.Item(3, 8).Value = 4.3  # H3
$ws.Cells.Item(3, 9).Value = 5  # I3
$ws.Cells.Item(3, 10).Value = 2.3  # J3
$ws.Cells.Item(3, 11).Value = 2.5  # K3
$ws.Cells.Item(3, 12).Value = 0  # L3
$ws.Cells.Item(3, 13).Value = 1.4  # M3
$ws.Cells.Item(3, 14).Value = 1.45  # N3
$ws.Cells.Item(3, 15).Value = 3  # O3
$ws.Cells.Item(3, 16).Value = 1.12  # P3
$ws.Cells.Item(3, 17).Value = 8.2  # Q3
$ws.Cells.Item(3, 18).Value = 1.03  # R3
$ws.Cells.Item(3, 19).Value = 32  # S3
$ws.Cells.Item(3, 20).Value = 4.3  # T3
$ws.Cells.Item(3, 21).Value = 1.24  # U3
$ws.Cells.Item(3, 22).Value = 1.23  # V3
$ws.Cells.Item(3, 23).Value = 1.52  # W3
$ws.Cells.Item(3, 24).Value = 3.65  # X3
$ws.Cells.Item(3, 25).Value = 8  # Y3
$ws.Cells.Item(3, 26).Value = 42  # Z3
$ws.Cells.Item(3, 27).Value = 320  # AA3
$ws.Cells.Item(3, 28).Value = 4.9  # AB3
$ws.Cells.Item(3, 29).Value = 9.6  # AC3
$ws.Cells.Item(3, 30).Value = 48  # AD3
$ws.Cells.Item(3, 31).Value = 350  # AE3
$ws.Cells.Item(3, 32).Value = 16  # AF3
$ws.Cells.Item(3, 33).Value = 34  # AG3
$ws.Cells.Item(3, 34).Value = 120  # AH3
$ws.Cells.Item(3, 35).Value = 1000  # AI3
$ws.Cells.Item(3, 36).Value = 90  # AJ3
$ws.Cells.Item(3, 37).Value = 170  # AK3
$ws.Cells.Item(3, 38).Value = 620  # AL3
$ws.Cells.Item(3, 39).Value = 1000  # AM3
$ws.Cells.Item(3, 40).Value = 310  # AN3
$ws.Cells.Item(3, 41).Value = 1000  # AO3

# Row 4
$ws.Cells.Item(4, 6).Value = 3.05  # F4
$ws.Cells.Item(4, 7).Value = 3.15  # G4
$ws.Cells.Item(4, 8).Value = 2.24  # H4
$ws.Cells.Item(4, 9).Value = 2.28  # I4
$ws.Cells.Item(4, 10).Value = 4.2  # J4
$ws.Cells.Item(4, 11).Value = 4.4  # K4
$ws.Cells.Item(4, 12).Value = 1.26  # L4
$ws.Cells.Item(4, 13).Value = 1.03  # M4
$ws.Cells.Item(4, 14).Value = 7.2  # N4
$ws.Cells.Item(4, 15).Value = 1.15  # O4
$ws.Cells.Item(4, 16).Value = 3.15  # P4
$ws.Cells.Item(4, 17).Value = 1.45  # Q4
$ws.Cells.Item(4, 18).Value = 1.85  # R4
$ws.Cells.Item(4, 19).Value = 2.16  # S4
$ws.Cells.Item(4, 20).Value = 1.46  # T4
$ws.Cells.Item(4, 21).Value = 3.1  # U4
$ws.Cells.Item(4, 22).Value = 1.79  # V4
$ws.Cells.Item(4, 23).Value = 1.46  # W4
$ws.Cells.Item(4, 24).Value = 30  # X4
$ws.Cells.Item(4, 25).Value = 18  # Y4
$ws.Cells.Item(4, 26).Value = 20  # Z4
$ws.Cells.Item(4, 27).Value = 30  # AA4
$ws.Cells.Item(4, 28).Value = 22  # AB4
$ws.Cells.Item(4, 29).Value = 10.5  # AC4
$ws.Cells.Item(4, 30).Value = 11.5  # AD4
$ws.Cells.Item(4, 31).Value = 19.5  # AE4
$ws.Cells.Item(4, 32).Value = 28  # AF4
$ws.Cells.Item(4, 33).Value = 14  # AG4
$ws.Cells.Item(4, 34).Value = 13.5  # AH4
$ws.Cells.Item(4, 35).Value = 23  # AI4
$ws.Cells.Item(4, 36).Value = 48  # AJ4
$ws.Cells.Item(4, 37).Value = 27  # AK4
$ws.Cells.Item(4, 38).Value = 29  # AL4
$ws.Cells.Item(4, 39).Value = 42  # AM4
$ws.Cells.Item(4, 40).Value = 14.5  # AN4
$ws.Cells.Item(4, 41).Value = 9  # AO4

